$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet used to carry a two-row header ("...mation"/"pompes)" style
# labels on row 1, units on row 2) above 13 data rows (rows 3-15). The new
# layout collapses this into a single header row followed directly by the
# 13 data rows (rows 2-14) -- the data rows themselves are unchanged, so
# removing the old row 2 (units row) shifts everything up by exactly one
# row and leaves the data rows already correct.
$ws.Rows(2).Delete()

# Rewrite row 1 with the new column headers.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# A1:E1 (the idx/idx2/Name/Date Start/Date End labels) use the plain
# default look (10pt, no special number format).
$ws.Range("A1:E1").ClearFormats()

# F1:K1 (the unit headers) use the smaller 9pt font applied to the rest
# of the data table.
$ws.Range("F1:K1").Font.Size = 9

# Match the new active selection left by the edit.
$ws.Range("A2:K2").Select()
